$d = $word.ActiveDocument
# p17 is right before p18 (numId=3). Let's check p17's ListFormat.
$p17 = $d.Paragraphs.Item(17)
Write-Output ("p17 text: " + $p17.Range.Text.Substring(0,30))
$p18 = $d.Paragraphs.Item(18)
Write-Output ("p18 ListID: " + $p18.Range.ListFormat.List.ListID)

$p10 = $d.Paragraphs.Item(10)
$lf10 = $p10.Range.ListFormat
# Try continue previous list referencing p18's template while p10 is NOT adjacent -- confirm fails (already know)
# Now let's check CanContinuePreviousList
Write-Output ("CanContinuePreviousList defined? ")
